$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 2884.8572
$ws.Range("I64").Value = 2499
$ws.Range("J64").Value = 3039.2
$ws.Range("K64").Value = 2499
$ws.Range("L64").Value = 3039.2
$ws.Range("M64").Value = -2251
$ws.Range("N64").Value = -3535.2

# Row 67
$ws.Range("H67").Value = 2884.8572
$ws.Range("I67").Value = 2499
$ws.Range("J67").Value = 3039.2
$ws.Range("K67").Value = 2499
$ws.Range("L67").Value = 3039.2
$ws.Range("M67").Value = -1641
$ws.Range("N67").Value = -4755.2

# Row 69
$ws.Range("H69").Value = 440396.2
$ws.Range("I69").Value = 631494.7
$ws.Range("K69").Value = 1894484.1
$ws.Range("M69").Value = -1893610.1

# Row 72
$ws.Range("H72").Value = 440396.2
$ws.Range("I72").Value = 631494.7
$ws.Range("K72").Value = 5683452.3
$ws.Range("M72").Value = -5679084.3

# Row 76
$ws.Range("H76").Value = 7711
$ws.Range("J76").Value = 6663.3335
$ws.Range("L76").Value = 6663.3335
$ws.Range("N76").Value = -7293.3335

# Row 79
$ws.Range("H79").Value = 7711
$ws.Range("J79").Value = 6663.3335
$ws.Range("L79").Value = 6663.3335
$ws.Range("N79").Value = -8847.333500000001

# Row 86
$ws.Range("H86").Value = 1500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1500
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -3746

# Row 89
$ws.Range("H89").Value = 1500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 7500
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -18732

# Row 98
$ws.Range("H98").Value = 54261.383
$ws.Range("I98").Value = 31357.268
$ws.Range("J98").Value = 226042.25
$ws.Range("K98").Value = 31357.268
$ws.Range("L98").Value = 226042.25
$ws.Range("M98").Value = -29859.268
$ws.Range("N98").Value = -229038.25

# Row 122
$ws.Range("H122").Value = 54261.383
$ws.Range("I122").Value = 31357.268
$ws.Range("J122").Value = 226042.25
$ws.Range("K122").Value = 94071.804
$ws.Range("L122").Value = 678126.75
$ws.Range("M122").Value = -91621.804
$ws.Range("N122").Value = -683026.75

# Row 138
$ws.Range("H138").Value = 3653.389
$ws.Range("I138").Value = 2117.25
$ws.Range("J138").Value = 4092.2856
$ws.Range("K138").Value = 6351.75
$ws.Range("L138").Value = 12276.8568
$ws.Range("M138").Value = -1211.75
$ws.Range("N138").Value = -22556.8568

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 19860.271
$ws.Range("I32").Value = 22333.5
$ws.Range("K32").Value = 22333.5
$ws.Range("M32").Value = -22046.5

# Row 74
$ws.Range("H74").Value = 1805.6744
$ws.Range("I74").Value = 1568.8422
$ws.Range("J74").Value = 3605.6
$ws.Range("K74").Value = 1568.8422
$ws.Range("L74").Value = 3605.6
$ws.Range("M74").Value = -694.8422
$ws.Range("N74").Value = -5353.6

# Row 77
$ws.Range("H77").Value = 1805.6744
$ws.Range("I77").Value = 1568.8422
$ws.Range("J77").Value = 3605.6
$ws.Range("K77").Value = 7844.211
$ws.Range("L77").Value = 18028
$ws.Range("M77").Value = -3476.211
$ws.Range("N77").Value = -26764

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2243.2856
$ws.Range("I105").Value = 1703.9
$ws.Range("K105").Value = 1703.9
$ws.Range("M105").Value = 43.09999999999991

# Row 134
$ws.Range("H134").Value = 3809.1667
$ws.Range("I134").Value = 3050.875
$ws.Range("J134").Value = 9875.5
$ws.Range("K134").Value = 9152.625
$ws.Range("L134").Value = 29626.5
$ws.Range("M134").Value = -6617.625
$ws.Range("N134").Value = -34696.5

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 66569.12
$ws.Range("I62").Value = 106237.1
$ws.Range("J62").Value = 9900.571
$ws.Range("K62").Value = 106237.1
$ws.Range("L62").Value = 9900.571
$ws.Range("M62").Value = -105613.1
$ws.Range("N62").Value = -11148.571

# Row 65
$ws.Range("H65").Value = 66569.12
$ws.Range("I65").Value = 106237.1
$ws.Range("J65").Value = 9900.571
$ws.Range("K65").Value = 531185.5
$ws.Range("L65").Value = 49502.855
$ws.Range("M65").Value = -528065.5
$ws.Range("N65").Value = -55742.855

# Row 107
$ws.Range("H107").Value = 1045.9697
$ws.Range("I107").Value = 594.95557
$ws.Range("K107").Value = 594.95557
$ws.Range("M107").Value = 1325.04443

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5429.5713
$ws.Range("I70").Value = 5552
$ws.Range("K70").Value = 5552
$ws.Range("M70").Value = -5282

# Row 73
$ws.Range("H73").Value = 5429.5713
$ws.Range("I73").Value = 5552
$ws.Range("K73").Value = 5552
$ws.Range("M73").Value = -4616

# Row 80
$ws.Range("H80").Value = 11276.056
$ws.Range("I80").Value = 4141
$ws.Range("J80").Value = 14020.308
$ws.Range("K80").Value = 4141
$ws.Range("L80").Value = 14020.308
$ws.Range("M80").Value = -3143
$ws.Range("N80").Value = -16016.308

# Row 83
$ws.Range("H83").Value = 11276.056
$ws.Range("I83").Value = 4141
$ws.Range("J83").Value = 14020.308
$ws.Range("K83").Value = 20705
$ws.Range("L83").Value = 70101.54000000001
$ws.Range("M83").Value = -15713
$ws.Range("N83").Value = -80085.54000000001

# Row 109
$ws.Range("H109").Value = 25070.072
$ws.Range("J109").Value = 25070.072
$ws.Range("L109").Value = 25070.072
$ws.Range("N109").Value = -27150.072

# Row 122
$ws.Range("H122").Value = 3285.5312
$ws.Range("I122").Value = 2592.238
$ws.Range("J122").Value = 4609.091
$ws.Range("K122").Value = 7776.714
$ws.Range("L122").Value = 13827.273
$ws.Range("M122").Value = -5326.714
$ws.Range("N122").Value = -18727.273

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 31251946
$ws.Range("I16").Value = 43479180
$ws.Range("K16").Value = 43479180
$ws.Range("M16").Value = -43479010

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# Row 81
$ws.Range("H81").Value = 4498.5713
$ws.Range("I81").Value = 1333.3334
$ws.Range("J81").Value = 6872.5
$ws.Range("K81").Value = 2666.6668
$ws.Range("L81").Value = 13745
$ws.Range("M81").Value = -1605.6668
$ws.Range("N81").Value = -15867

# Row 84
$ws.Range("H84").Value = 4498.5713
$ws.Range("I84").Value = 1333.3334
$ws.Range("J84").Value = 6872.5
$ws.Range("K84").Value = 13333.334
$ws.Range("L84").Value = 68725
$ws.Range("M84").Value = -8029.333999999999
$ws.Range("N84").Value = -79333

# Row 107
$ws.Range("H107").Value = 330.96667
$ws.Range("I107").Value = 251.84616
$ws.Range("J107").Value = 845.25
$ws.Range("K107").Value = 755.5384799999999
$ws.Range("L107").Value = 2535.75
$ws.Range("M107").Value = 1164.46152
$ws.Range("N107").Value = -6375.75
